$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").Value = '26.386.64'
$ws.Range("E2").Value = '  -0.48%  '

# Row 3: update D3, E3
$ws.Range("D3").Value = '1.831.94'
$ws.Range("E3").Value = '  -0.40%  '

# Row 4: update D4, E4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.07%  '

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.01'
$ws.Range("E5").Value = '  -3.21%  '

# Row 6: update E6
$ws.Range("E6").Value = '  +0.12%  '

# Row 7: update D7, E7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5255'
$ws.Range("E7").Value = '  -0.09%  '

# Row 8: update D8, E8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2783'
$ws.Range("E8").Value = '  -13.24%  '

# Row 9: update D9, E9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06843'
$ws.Range("E9").Value = '  +0.81%  '

# Row 10: update D10, E10
$ws.Range("D10").Value = '1.852.63'
$ws.Range("E10").Value = '  +0.67%  '

# Row 11: update D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.53'
$ws.Range("E11").Value = '  -11.86%  '

# Row 12: update D12, E12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07087'
$ws.Range("E12").Value = '  -8.39%  '

# Row 13: update D13, E13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6882'
$ws.Range("E13").Value = '  -12.25%  '

# Row 14: update D14, E14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '86.07'
$ws.Range("E14").Value = '  -1.88%  '

# Row 15: update D15, E15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.849'
$ws.Range("E15").Value = '  -3.19%  '

# Row 16: update D16, E16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.08%  '

# Row 17: update E17
$ws.Range("E17").Value = '  +0.06%  '

# Row 18: update D18, E18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.17'
$ws.Range("E18").Value = '  -4.78%  '

# Row 19: update D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007326'
$ws.Range("E19").Value = '  -7.52%  '

# Row 20: update D20, E20
$ws.Range("D20").Value = '26.420.81'
$ws.Range("E20").Value = '  -0.43%  '

# Row 21: update D21, E21
$ws.Range("D21").Value = '2.084.73'
$ws.Range("E21").Value = '  +0.33%  '

# Row 22: update D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.492'
$ws.Range("E22").Value = '  -2.73%  '

# Row 23: update D23, E23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.802'
$ws.Range("E23").Value = '  -2.90%  '

# Row 24: update D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.953'
$ws.Range("E24").Value = '  -4.48%  '

# Row 25: update D25, E25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.41'
$ws.Range("E25").Value = '  +0.73%  '

# Row 26: update D26, E26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.673'
$ws.Range("E26").Value = '  -0.48%  '

# Row 27: update D27, E27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.020'
$ws.Range("E27").Value = '  -6.19%  '

# Row 28: update D28, E28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.52'
$ws.Range("E28").Value = '  -2.35%  '

# Row 29: update D29, E29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '109.04'
$ws.Range("E29").Value = '  -2.20%  '

# Row 30: update D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.057'
$ws.Range("E30").Value = '  -2.15%  '

# Row 31: update D31, E31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08716'
$ws.Range("E31").Value = '  +0.30%  '

# Row 32: update D32, E32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.851'
$ws.Range("E32").Value = '  -5.26%  '

# Row 33: update D33, E33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04677'
$ws.Range("E33").Value = '  -3.82%  '

# Row 34: update D34, E34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.881'
$ws.Range("E34").Value = '  +0.93%  '

# Row 35: update B35, C35, D35, E35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.104'
$ws.Range("E35").Value = '  -2.57%  '

# Row 36: update B36, C36, D36, E36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7016'
$ws.Range("E36").Value = '  -3.91%  '

# Row 37: update D37, E37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.051'
$ws.Range("E37").Value = '  -1.15%  '

# Row 38: update D38, E38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.172'
$ws.Range("E38").Value = '  -3.26%  '

# Row 39: update D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01634'
$ws.Range("E39").Value = '  -7.16%  '

# Row 40: update D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4449'
$ws.Range("E40").Value = '  -6.94%  '

# Row 41: update D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8589'
$ws.Range("E41").Value = '  -3.64%  '

# Row 42: update D42, E42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '104.99'
$ws.Range("E42").Value = '  -4.30%  '

# Row 43: update D43, E43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.02%  '

# Row 44: update D44, E44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.728'
$ws.Range("E44").Value = '  -3.20%  '

# Row 45: update D45, E45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.984'
$ws.Range("E45").Value = '  -8.75%  '

# Row 46: update D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.570'
$ws.Range("E46").Value = '  -4.15%  '

# Row 47: update D47, E47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05569'
$ws.Range("E47").Value = '  -4.81%  '

# Row 48: update D48, E48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '58.55'
$ws.Range("E48").Value = '  -1.69%  '

# Row 49: update D49, E49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.31'
$ws.Range("E49").Value = '  -4.34%  '

# Row 50: update B50, C50, D50, E50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1173'
$ws.Range("E50").Value = '  -4.73%  '

# Row 51: update B51, C51, D51, E51
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8592'
$ws.Range("E51").Value = '  -3.99%  '
